# Update metadata sheet: switch the "horas-trabajadas" (col A) and
# "ocupacion-1-digito-descripcion" (col C) columns from measure/xsd:string
# to dimension/skos:Concept, and register their new mapping files in a
# freshly added row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: iaest-measure:... -> iaest-dimension:... for columns A and C
$ws.Range("A3").Value = "iaest-dimension:horas-trabajadas"
$ws.Range("C3").Value = "iaest-dimension:ocupacion-1-digito-descripcion"

# Row 4: "medida" -> "dim" for columns A and C (column E already says "dim")
$ws.Range("A4").Value = "dim"
$ws.Range("C4").Value = "dim"

# Row 5: "xsd:string" -> "skos:Concept" for columns A and C
$ws.Range("A5").Value = "skos:Concept"
$ws.Range("C5").Value = "skos:Concept"

# Row 6 (new): reuse the existing cell formatting, then fill in the new
# mapping-file references (only columns A and C get a value, as in the
# source report).
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("A6").Value = "mapping-horas-trabajadas.xlsx"
$ws.Range("C6").Value = "mapping-ocupacion-1-digito-descripcion.xlsx"
